# Generate Report for Handoff
#
# The "b.md" file (in both the zh-cn and de-de localization targets) has just
# gone through a new handoff. Update the Overview sheet and the per-locale
# status sheets to reflect the new "Ready for handoff" status, the freshly
# generated handoff xliff files / timestamps, and the stale-handback error
# that was detected for that file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the b.md file. Its zh-cn / de-de status columns
# (E, F) move from "Handed back: in sync with en-US" to "Ready for handoff",
# and the "Latest HO Xliff Generate Date" column (G) gets the new timestamp.
# ---------------------------------------------------------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-21 16:45:47"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) picks up the new handoff file/timestamp, the
# refreshed status, the no-longer-duplicate content flag, and the newly
# detected stale-handback error message.
# ---------------------------------------------------------------------------
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 16:45:43"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b564cfe5d8fdf6ee07b2623916e877437318fc7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81a0683631d0215932090c41a3d4e84b9247fdf3/e2e/b.md."

# Widen the Error Detail column so the long message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets the same treatment as zh-cn above.
# ---------------------------------------------------------------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 16:45:47"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b564cfe5d8fdf6ee07b2623916e877437318fc7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81a0683631d0215932090c41a3d4e84b9247fdf3/e2e/b.md."

# Widen the Error Detail column so the long message is readable.
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
